$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 152 (id 150): shift in data that was previously for the match
# that is now row 154 in the old layout, refreshed with new odds.
$ws.Range("B152").Value = 6965778
$ws.Range("E152").Value = 45403.5
$ws.Range("F152").Value = "NK Lokomotiva Zagreb"
$ws.Range("G152").Value = "Dinamo Zagreb"
$ws.Range("K152").Value = 6
$ws.Range("L152").Value = 4.333
$ws.Range("M152").Value = 1.5
$ws.Range("N152").Value = 7
$ws.Range("O152").Value = 3.8
$ws.Range("P152").Value = 1.5
$ws.Range("Q152").Value = 1
$ws.Range("R152").Value = 2.05
$ws.Range("S152").Value = 1.8
$ws.Range("T152").Value = 2.5
$ws.Range("U152").Value = 1.825
$ws.Range("V152").Value = 2.025

# Update row 153 (id 151)
$ws.Range("B153").Value = 6962506
$ws.Range("E153").Value = 45403.60416666666
$ws.Range("F153").Value = "HNK Rijeka"
$ws.Range("G153").Value = "HNK Gorica"
$ws.Range("K153").Value = 1.285
$ws.Range("L153").Value = 6
$ws.Range("M153").Value = 8.5
$ws.Range("N153").Value = 1.2
$ws.Range("O153").Value = 7
$ws.Range("P153").Value = 11
$ws.Range("Q153").Value = -1.75
$ws.Range("R153").Value = 1.8
$ws.Range("S153").Value = 2.05
$ws.Range("T153").Value = 3
$ws.Range("U153").Value = 1.975
$ws.Range("V153").Value = 1.875

# Update row 154 (id 152)
$ws.Range("B154").Value = 6957866
$ws.Range("E154").Value = 45404.54166666666
$ws.Range("F154").Value = "NK Rudes"
$ws.Range("G154").Value = "NK Varazdin"
$ws.Range("K154").Value = 5
$ws.Range("L154").Value = 3.75
$ws.Range("M154").Value = 1.615
$ws.Range("N154").Value = 4
$ws.Range("O154").Value = 3.6
$ws.Range("P154").Value = 1.85
$ws.Range("Q154").Value = 0.5
$ws.Range("R154").Value = 1.925
$ws.Range("S154").Value = 1.925
$ws.Range("T154").Value = 2.5
$ws.Range("U154").Value = 1.95
$ws.Range("V154").Value = 1.9

# Update row 155 (id 153) - this becomes a brand new fixture
$ws.Range("B155").Value = 6990513
$ws.Range("E155").Value = 45408.54166666666
$ws.Range("F155").Value = "Istra 1961"
$ws.Range("G155").Value = "Slaven Belupo"
$ws.Range("K155").Value = 2.1
$ws.Range("L155").Value = 3.1
$ws.Range("M155").Value = 3.8
$ws.Range("N155").Value = 2.1
$ws.Range("O155").Value = 3.1
$ws.Range("P155").Value = 3.8
$ws.Range("Q155").Value = -0.25
$ws.Range("R155").Value = 1.8
$ws.Range("S155").Value = 2.05
$ws.Range("T155").Value = 2.25
$ws.Range("U155").Value = 2.025
$ws.Range("V155").Value = 1.825

# Row 156 (the old last fixture) is no longer present - remove it entirely,
# which also shifts the sheet dimension from AC156 down to AC155.
$ws.Rows("156:156").Delete()
